$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet (3rd sheet) ---
# Insert a new (blank) column at N, shifting the old N/O/P columns to O/P/Q.
$wsRepay = $wb.Worksheets.Item(3)
$wsRepay.Columns("N:N").Insert()

# Match the inserted column's width to its left neighbour (column M), as Excel does.
$wsRepay.Columns("N:N").ColumnWidth = $wsRepay.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with R13 selected.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("R13").Select() | Out-Null

# --- "Edit Repayment Schedule" sheet (5th sheet) ---
# It is no longer the active tab; update its remembered selection to C15.
$wsEdit = $wb.Worksheets.Item(5)
$wsEdit.Activate() | Out-Null
$wsEdit.Range("C15").Select() | Out-Null

# Restore "Repayment schedule" as the active sheet so it is the one shown/selected.
$wsRepay.Activate() | Out-Null
